$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Prepare new row 13 so it inherits row 12's current formatting ---
# (row 12 already carries the "data row" style used by A/B/C = s6, D = s4,
#  with ht=32.1/customHeight, and that is exactly what the new row 13 needs).
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("13:13").RowHeight = $ws.Rows("12:12").RowHeight
$excel.CutCopyMode = 0

# --- Move the old row 12 content ("I MIND MAP") down into row 13 ---
$ws.Range("A13").Value = $ws.Range("A12").Value2

# --- Clear the old row 11 content ("LUCID CHART" / "Grátis e Pago") ---
# it is being relocated to row 12 below.
$ws.Range("A11").Value = $null
$ws.Range("B11").Value = $null

# --- Fill row 12 with the relocated LUCID CHART data plus the two new
#     columns (web / Somente On-line. Pagamento mensal) ---
$ws.Range("A12").Value = "LUCID CHART"
$ws.Range("B12").Value = "Grátis e Pago"
$ws.Range("C12").Value = "web"
$ws.Range("D12").Value = "Somente On-line. Pagamento mensal"

# --- Fill in the newly-visible "Observação" for the XMIND row ---
$ws.Range("D9").Value = "Pagamento mensal"

# --- Finish filling row 13 (I MIND MAP) with the new columns ---
$ws.Range("B13").Value = "Gratis e Pago"
$ws.Range("C13").Value = "Windows - Mac"
$ws.Range("D13").Value = "Pagamento Único em USD"

# --- Update the view state to match: scrolled so row 4 is at the top,
#     and the newly added row 13 selected ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A13:XFD13").Select()
